$wb = $excel.ActiveWorkbook

# --- Hoja1: update the conversion summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$text = $cellA1.Value()
$text = $text.Replace("1000 Bs = 12.96 = 52986.14 pesos", "1000 Bs = 12.97 = 52658.88 pesos")
$text = $text.Replace("52986.14 pesos = 12.97 = 980.73 Bs", "52658.88 pesos = 12.91 = 965.66 Bs")
$cellA1.Value = $text

# --- tasas: update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 77.09999999999999
$ws2.Range("O10").Value = 4060
$ws2.Range("N12").Value = 4079
$ws2.Range("O12").Value = 74.801
